$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = "ELT-1A-Desenho Técnico"
$ws.Range("E3").Value = "MEC-2A-CAD"
$ws.Range("F3").Value = "MEC-2A-CAD"

# Row 4
$ws.Range("E4").Value = "-"

# Row 6
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "MEC-1A-Desenho Técnico"

# Row 7
$ws.Range("B7").Value = "-"
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "MEC-1A-Desenho Técnico"
